# C24239 and C24240 scripts added
#
# "General" sheet (sheet1): two new automation rows (52 & 53) for the
# newly added test cases, each with a mailto hyperlink in column E
# (matching the existing rows' pattern).
#
# "Data" sheet (sheet2): one new row (19) holding the account number used
# by the C24240 test case.
#
# The Data sheet also becomes the active tab/selection, while General's
# selection moves to the newly appended A53 cell.

$wb = $excel.ActiveWorkbook

# ---- General sheet ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("General")

$ws1.Range("A52").Value = "C24239_VerifyAllAccountsAddressChange"
$ws1.Range("B52").Value = "Yes"
$ws1.Range("C52").Value = "Android"
$ws1.Range("D52").Value = "user2046417"
$ws1.Range("E52").Value = "Kony@1234"
$ws1.Hyperlinks.Add($ws1.Range("E52"), "mailto:Kony@1234") | Out-Null
$ws1.Range("E52").Style = "Hyperlink"

$ws1.Range("A53").Value = "C24240_VerifyAddressChangeSuccessMessage"
$ws1.Range("B53").Value = "Yes"
$ws1.Range("C53").Value = "Android"
$ws1.Range("D53").Value = "user2046417"
$ws1.Range("E53").Value = "Kony@1234"
$ws1.Hyperlinks.Add($ws1.Range("E53"), "mailto:Kony@1234") | Out-Null
$ws1.Range("E53").Style = "Hyperlink"

$ws1.Range("A53").Select() | Out-Null

# ---- Data sheet --------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Data")

$ws2.Range("A19").Value = "C24240_VerifyAddressChangeSuccessMessage"
$ws2.Range("B19").Value = "'20464178"
$ws2.Range("B19").WrapText = $true

$ws2.Activate() | Out-Null
$ws2.Range("B19").Select() | Out-Null
